$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("Q2").Value = 1.62
$ws.Range("AA6").Value = 19
$ws.Range("AH6").Value = 9.5
$ws.Range("AO6").Value = 11
$ws.Range("AP6").Value = 26
$ws.Range("AR6").Value = 67
$ws.Range("G6").Value = 1.91
$ws.Range("H6").Value = 3.1
$ws.Range("I6").Value = 4.75
$ws.Range("K6").Value = 1.95
$ws.Range("X6").Value = 7.5
$ws.Range("Z6").Value = 15
$ws.Range("G8").Value = 2.15
$ws.Range("N8").Value = 7.5
$ws.Range("AH9").Value = 8
$ws.Range("AI9").Value = 12
$ws.Range("AJ9").Value = 10
$ws.Range("AX9").Value = 15
$ws.Range("AZ9").Value = 51
$ws.Range("BC9").Value = 126
$ws.Range("G9").Value = 2.75
$ws.Range("I9").Value = 2.45
$ws.Range("J9").Value = 3.5
$ws.Range("L9").Value = 3.2
$ws.Range("N9").Value = 9.5
$ws.Range("Q11").Value = 1.9
$ws.Range("R11").Value = 1.95
$ws.Range("AC12").Value = 13
$ws.Range("AO12").Value = 7.5
$ws.Range("AP12").Value = 17
$ws.Range("AU12").Value = 8
$ws.Range("BA12").Value = 101
$ws.Range("K12").Value = 2.37
$ws.Range("U12").Value = 1.8
$ws.Range("V12").Value = 1.91
$ws.Range("G13").Value = 1.92
$ws.Range("Q13").Value = 1.84
$ws.Range("R13").Value = 1.89
$ws.Range("U17").Value = 1.92
$ws.Range("V17").Value = 1.77
$ws.Range("Q18").Value = 1.77
$ws.Range("R18").Value = 1.97
$ws.Range("U18").Value = 1.63
$ws.Range("Q19").Value = 1.69
$ws.Range("R19").Value = 2.07
$ws.Range("U19").Value = 1.58
$ws.Range("J20").Value = 2.37
$ws.Range("Q20").Value = 1.77
$ws.Range("U20").Value = 1.69
$ws.Range("K22").Value = 2.37
$ws.Range("O22").Value = 1.18
$ws.Range("P22").Value = 4.5
$ws.Range("Q22").Value = 1.65
$ws.Range("R22").Value = 2.2
$ws.Range("Q23").Value = 1.6
$ws.Range("R23").Value = 2.3
$ws.Range("AA25").Value = 19
$ws.Range("AC25").Value = 8
$ws.Range("AG25").Value = 351
$ws.Range("AN25").Value = 4
$ws.Range("AT25").Value = 2.5
$ws.Range("BA25").Value = 101
$ws.Range("BB25").Value = 251
$ws.Range("G25").Value = 2.25
$ws.Range("I25").Value = 3.4
$ws.Range("L25").Value = 4
$ws.Range("S25").Value = 1.5
$ws.Range("T25").Value = 2.5
$ws.Range("AC26").Value = 7
$ws.Range("AG26").Value = 501
$ws.Range("AI26").Value = 15
$ws.Range("AU26").Value = 9
$ws.Range("AW26").Value = 5
$ws.Range("AY26").Value = 34
$ws.Range("BB26").Value = 301
$ws.Range("G26").Value = 2.4
$ws.Range("J26").Value = 3.2
$ws.Range("L26").Value = 4
$ws.Range("M26").Value = 1.1
$ws.Range("N26").Value = 7
$ws.Range("Q26").Value = 2.5
$ws.Range("R26").Value = 1.5
$ws.Range("U26").Value = 2.05
$ws.Range("V26").Value = 1.7
$ws.Range("X26").Value = 10
$ws.Range("AH27").Value = 8
$ws.Range("AI27").Value = 13
$ws.Range("AK27").Value = 29
$ws.Range("AO27").Value = 15
$ws.Range("AX27").Value = 17
$ws.Range("G27").Value = 2.75
$ws.Range("I27").Value = 2.8
$ws.Range("J27").Value = 3.4
$ws.Range("L27").Value = 3.5
$ws.Range("O27").Value = 1.36
$ws.Range("P27").Value = 3
$ws.Range("Q27").Value = 2.25
$ws.Range("R27").Value = 1.62
$ws.Range("W27").Value = 7.5
$ws.Range("X27").Value = 12
$ws.Range("Z27").Value = 26
$ws.Range("Q31").Value = 2.08
$ws.Range("R31").Value = 1.73
$ws.Range("Q36").Value = 1.4
$ws.Range("R36").Value = 2.88
$ws.Range("AA37").Value = 17
$ws.Range("AB37").Value = 41
$ws.Range("AC37").Value = 6.5
$ws.Range("AH37").Value = 12
$ws.Range("AI37").Value = 29
$ws.Range("AO37").Value = 9
$ws.Range("AP37").Value = 26
$ws.Range("AQ37").Value = 34
$ws.Range("AR37").Value = 67
$ws.Range("AS37").Value = 251
$ws.Range("AT37").Value = 2.38
$ws.Range("AW37").Value = 7
$ws.Range("BC37").Value = 126
$ws.Range("G37").Value = 1.67
$ws.Range("H37").Value = 3.4
$ws.Range("I37").Value = 6.25
$ws.Range("J37").Value = 2.38
$ws.Range("K37").Value = 2
$ws.Range("M37").Value = 1.11
$ws.Range("N37").Value = 6.5
$ws.Range("O37").Value = 1.44
$ws.Range("P37").Value = 2.63
$ws.Range("Q37").Value = 2.4
$ws.Range("R37").Value = 1.53
$ws.Range("S37").Value = 1.53
$ws.Range("T37").Value = 2.38
$ws.Range("W37").Value = 5
$ws.Range("Z37").Value = 12
$ws.Range("AL39").Value = 29
$ws.Range("AW39").Value = 4.75
$ws.Range("AY39").Value = 29
$ws.Range("BD39").Value = 151
$ws.Range("G39").Value = 2.7
$ws.Range("H39").Value = 2.5
$ws.Range("I39").Value = 3.25
$ws.Range("M39").Value = 1.13
$ws.Range("N39").Value = 6
$ws.Range("Q39").Value = 2.6
$ws.Range("R39").Value = 1.48
$ws.Range("S39").Value = 1.5
$ws.Range("T39").Value = 2.37
$ws.Range("U39").Value = 1.91
$ws.Range("V39").Value = 1.8
$ws.Range("W39").Value = 7.5
$ws.Range("X39").Value = 12
$ws.Range("Z39").Value = 26
$ws.Range("S40").Value = 1.37
$ws.Range("S41").Value = 1.37
$ws.Range("N42").Value = 8
$ws.Range("Q42").Value = 2.3
$ws.Range("R42").Value = 1.6
$ws.Range("AI45").Value = 23
$ws.Range("AK45").Value = 41
$ws.Range("AN45").Value = 4
$ws.Range("AO45").Value = 9.5
$ws.Range("AQ45").Value = 29
$ws.Range("AV45").Value = 41
$ws.Range("AW45").Value = 6
$ws.Range("AX45").Value = 21
$ws.Range("G45").Value = 1.8
$ws.Range("H45").Value = 3.6
$ws.Range("I45").Value = 4.33
$ws.Range("J45").Value = 2.4
$ws.Range("K45").Value = 2.3
$ws.Range("L45").Value = 4.33
$ws.Range("Q45").Value = 1.7
$ws.Range("R45").Value = 2.1
$ws.Range("U45").Value = 1.62
$ws.Range("V45").Value = 2.2
$ws.Range("W45").Value = 9
$ws.Range("X45").Value = 10
$ws.Range("Z45").Value = 15
$ws.Range("AA46").Value = 23
$ws.Range("AD46").Value = 7.5
$ws.Range("AE46").Value = 12
$ws.Range("AG46").Value = 126
$ws.Range("AH46").Value = 10
$ws.Range("AK46").Value = 19
$ws.Range("AS46").Value = 126
$ws.Range("AT46").Value = 3.4
$ws.Range("AU46").Value = 7
$ws.Range("AY46").Value = 17
$ws.Range("BA46").Value = 41
$ws.Range("BC46").Value = 401
$ws.Range("G46").Value = 3.4
$ws.Range("H46").Value = 3.8
$ws.Range("J46").Value = 3.75
$ws.Range("K46").Value = 2.38
$ws.Range("L46").Value = 2.6
$ws.Range("O46").Value = 1.18
$ws.Range("P46").Value = 4.5
$ws.Range("Q46").Value = 1.62
$ws.Range("R46").Value = 2.25
$ws.Range("S46").Value = 1.27
$ws.Range("T46").Value = 3.4
$ws.Range("U46").Value = 1.57
$ws.Range("V46").Value = 2.25
$ws.Range("AE47").Value = 14.5
$ws.Range("AO47").Value = 8.75
$ws.Range("AQ47").Value = 30
$ws.Range("AT47").Value = 2.82
$ws.Range("AY47").Value = 29
$ws.Range("I47").Value = 4.05
$ws.Range("S47").Value = 1.38
$ws.Range("T47").Value = 2.82
$ws.Range("U47").Value = 1.72
$ws.Range("V47").Value = 2
$ws.Range("W47").Value = 7.6
$ws.Range("Z47").Value = 14
